# ISAICP-3188: Migrate documents with unlimited cardinality file field.
#
# Appends two new rows (30 & 31) to the "1. Content items" sheet describing
# a "Document with URL and multiple files" migration case, grows the
# Table18915 Excel table / autofilter / hidden _FilterDatabase defined name
# to cover the new rows, wires up the new collection-owner mailto
# hyperlink on M30, and leaves the selection on A2 (as seen in the
# post-edit workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows ------------------------------------------------------
# Write the row 31 text first so the new shared strings are appended in
# the same order as the reference edit (" " before the longer string).
$ws.Cells.Item(31, 3).Value = " "

$ws.Cells.Item(30, 1).Value = "Document"
$ws.Cells.Item(30, 2).Value = 125548
$ws.Cells.Item(30, 3).Value = "Document with URL and multiple files"
$ws.Cells.Item(30, 4).Value = "Archived collection"
$ws.Cells.Item(30, 6).Value = "Open government"
$ws.Cells.Item(30, 7).Value = "No"
$ws.Cells.Item(30, 8).Value = "Yes"
$ws.Cells.Item(30, 13).Value = "doe@example.com"

# --- Hyperlink for the new Collection Owner cell -------------------------
$null = $ws.Hyperlinks.Add($ws.Range("M30"), "mailto:doe@example.com")

# --- Grow the table / autofilter to include the new rows -----------------
$lo = $ws.ListObjects.Item("Table18915")
$lo.Resize($ws.Range("A1:Q31"))

# --- Keep the hidden AutoFilter defined name in sync ----------------------
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "='1. Content items'!`$A`$1:`$Q`$31"

# --- Restore the selection recorded in the edited workbook ----------------
$null = $ws.Range("A2").Select()
